$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96; this shifts existing rows 96..129 down to 97..130
$ws.Rows.Item(96).Insert()

# Populate the new row 96 with the new weekly price record
$ws.Range("A96").Value = 7
$ws.Range("B96").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C96").Value = "Ñuble"
$ws.Range("D96").Value = 44463
$ws.Range("E96").Value = 16
$ws.Range("F96").Value = 100112032
$ws.Range("G96").Value = "Zapallo italiano"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 120
$ws.Range("K96").Value = 14000
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = 14500
$ws.Range("N96").Value = "$/caja 50 unidades"
$ws.Range("O96").Value = "Región de Arica y Parinacota"
$ws.Range("P96").Value = 290
$ws.Range("Q96").Value = 50
$ws.Range("R96").Value = "Hortaliza"
